# resources/Files/ExcelReader.xlsx -- "add all new file"
#
# Fix the "Passwrod" header typo, drop the stray leftover number in E2,
# add a new "imat0rik" password sample row, and let column C re-fit its
# (now shorter) contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("authentication")

# C1: "Passwrod" -> "Password" (typo fix)
$ws.Range("C1").Value = "Password"

# Remove the stray "45" that used to live in E2
$ws.Range("E2").ClearContents()

# New sample row: C3 = "imat0rik"
$ws.Range("C3").Value = "imat0rik"

# Column C no longer needs to be as wide (values got shorter) - autofit it
$ws.Columns.Item(3).AutoFit()

# Move the active selection to C1
$ws.Range("C1").Select()
